$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: 23-03-2025, Sunrisers Hyderabad vs Rajasthan Royals, toss=Rajasthan Royals, winner=Rajasthan Royals, login=vijay
$ws.Range("A3").Value = "23-03-2025"
$ws.Range("B3").Value = "Sunrisers Hyderabad vs Rajasthan Royals"
$ws.Range("C3").Value = "Rajasthan Royals"
$ws.Range("D3").Value = "Rajasthan Royals"
$ws.Range("E3").Value = "vijay"

# Row 4: 23-03-2025, Chennai Super Kings vs Mumbai Indians, toss=Mumbai Indians, winner=Mumbai Indians, login=vijay
$ws.Range("A4").Value = "23-03-2025"
$ws.Range("B4").Value = "Chennai Super Kings vs Mumbai Indians"
$ws.Range("C4").Value = "Mumbai Indians"
$ws.Range("D4").Value = "Mumbai Indians"
$ws.Range("E4").Value = "vijay"

# Remove old row 5 (Chennai Super Kings vs Mumbai Indians, nandini) entirely
$ws.Rows.Item(5).Delete()
